# Remove the "ddays(...) / dyears(...)" table row (and its
# "date1 + ddays(...)" example cell) from the dates/durations table.
# The preceding row ("days(...)/years(...)" -> "date1 + days(...)")
# and the following row ("interval(...)") are left untouched.

$d = $word.ActiveDocument
$table = $d.Tables.Item(1)

for ($i = $table.Rows.Count; $i -ge 1; $i--) {
    $row = $table.Rows.Item($i)
    if ($row.Cells.Item(1).Range.Text -like "*ddays(...)*") {
        $row.Delete()
    }
}
